$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 150.0354306666667
$ws.Range("H2").Value = 450.106292
$ws.Range("I2").Value = 0.4152507364956075
$ws.Range("J2").Value = 0.4152507364956075
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 14.561928
$ws.Range("N2").Value = 43.685784
$ws.Range("O2").Value = 0.3501777048818433
$ws.Range("P2").Value = 0.3501777048818433
$ws.Range("Q2").Value = 2184.805138816992
$ws.Range("R2").Value = 19663.24624935293
$ws.Range("S2").Value = 0.1454115498565269
$ws.Range("T2").Value = 0.1454115498565269

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 150.0354306666667
$ws.Range("H3").Value = 450.106292
$ws.Range("I3").Value = 0.4152507364956075
$ws.Range("J3").Value = 0.4152507364956075
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 14.40015733333333
$ws.Range("N3").Value = 43.200472
$ws.Range("O3").Value = 0.3462875276490937
$ws.Range("P3").Value = 0.3462875276490937
$ws.Range("Q3").Value = 2160.533807174425
$ws.Range("R3").Value = 19444.80426456983
$ws.Range("S3").Value = 0.1437961508955292
$ws.Range("T3").Value = 0.1437961508955292

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 150.0354306666667
$ws.Range("H4").Value = 450.106292
$ws.Range("I4").Value = 0.4152507364956075
$ws.Range("J4").Value = 0.4152507364956075
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.886742333333333
$ws.Range("N4").Value = 8.660226999999999
$ws.Range("O4").Value = 0.06941888497676431
$ws.Range("P4").Value = 0.06941888497676431
$ws.Range("Q4").Value = 433.1136292053648
$ws.Range("R4").Value = 3898.022662848284
$ws.Range("S4").Value = 0.02882624311330524
$ws.Range("T4").Value = 0.02882624311330524

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 150.0354306666667
$ws.Range("H5").Value = 450.106292
$ws.Range("I5").Value = 0.4152507364956075
$ws.Range("J5").Value = 0.4152507364956075
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 9.735567333333334
$ws.Range("N5").Value = 29.206702
$ws.Range("O5").Value = 0.2341158824922987
$ws.Range("P5").Value = 0.2341158824922987
$ws.Range("Q5").Value = 1460.680037640998
$ws.Range("R5").Value = 13146.12033876898
$ws.Range("S5").Value = 0.09721679263024614
$ws.Range("T5").Value = 0.09721679263024613

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 68.382243
$ws.Range("H6").Value = 205.146729
$ws.Range("I6").Value = 0.1892604742946246
$ws.Range("J6").Value = 0.1892604742946246
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 14.561928
$ws.Range("N6").Value = 43.685784
$ws.Range("O6").Value = 0.3501777048818433
$ws.Range("P6").Value = 0.3501777048818433
$ws.Range("Q6").Value = 995.777299044504
$ws.Range("R6").Value = 8961.995691400536
$ws.Range("S6").Value = 0.06627479851334075
$ws.Range("T6").Value = 0.06627479851334074

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 68.382243
$ws.Range("H7").Value = 205.146729
$ws.Range("I7").Value = 0.1892604742946246
$ws.Range("J7").Value = 0.1892604742946246
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 14.40015733333333
$ws.Range("N7").Value = 43.200472
$ws.Range("O7").Value = 0.3462875276490937
$ws.Range("P7").Value = 0.3462875276490937
$ws.Range("Q7").Value = 984.7150580062321
$ws.Range("R7").Value = 8862.435522056088
$ws.Range("S7").Value = 0.06553854172518042
$ws.Range("T7").Value = 0.06553854172518039

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 68.382243
$ws.Range("H8").Value = 205.146729
$ws.Range("I8").Value = 0.1892604742946246
$ws.Range("J8").Value = 0.1892604742946246
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.886742333333333
$ws.Range("N8").Value = 8.660226999999999
$ws.Range("O8").Value = 0.06941888497676431
$ws.Range("P8").Value = 0.06941888497676431
$ws.Range("Q8").Value = 197.401915716387
$ws.Range("R8").Value = 1776.617241447483
$ws.Range("S8").Value = 0.01313825109570641
$ws.Range("T8").Value = 0.0131382510957064

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 68.382243
$ws.Range("H9").Value = 205.146729
$ws.Range("I9").Value = 0.1892604742946246
$ws.Range("J9").Value = 0.1892604742946246
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 9.735567333333334
$ws.Range("N9").Value = 29.206702
$ws.Range("O9").Value = 0.2341158824922987
$ws.Range("P9").Value = 0.2341158824922987
$ws.Range("Q9").Value = 665.739931130862
$ws.Range("R9").Value = 5991.659380177758
$ws.Range("S9").Value = 0.04430888296039706
$ws.Range("T9").Value = 0.04430888296039705

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 104.737245
$ws.Range("H10").Value = 314.211735
$ws.Range("I10").Value = 0.2898796499701289
$ws.Range("J10").Value = 0.2898796499701289
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 14.561928
$ws.Range("N10").Value = 43.685784
$ws.Range("O10").Value = 0.3501777048818433
$ws.Range("P10").Value = 0.3501777048818433
$ws.Range("Q10").Value = 1525.17622060836
$ws.Range("R10").Value = 13726.58598547524
$ws.Range("S10").Value = 0.1015093905184919
$ws.Range("T10").Value = 0.1015093905184918

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 104.737245
$ws.Range("H11").Value = 314.211735
$ws.Range("I11").Value = 0.2898796499701289
$ws.Range("J11").Value = 0.2898796499701289
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 14.40015733333333
$ws.Range("N11").Value = 43.200472
$ws.Range("O11").Value = 0.3462875276490937
$ws.Range("P11").Value = 0.3462875276490937
$ws.Range("Q11").Value = 1508.23280665988
$ws.Range("R11").Value = 13574.09525993892
$ws.Range("S11").Value = 0.1003817073039406
$ws.Range("T11").Value = 0.1003817073039406

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 104.737245
$ws.Range("H12").Value = 314.211735
$ws.Range("I12").Value = 0.2898796499701289
$ws.Range("J12").Value = 0.2898796499701289
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 2.886742333333333
$ws.Range("N12").Value = 8.660226999999999
$ws.Range("O12").Value = 0.06941888497676431
$ws.Range("P12").Value = 0.06941888497676431
$ws.Range("Q12").Value = 302.3494390182049
$ws.Range("R12").Value = 2721.144951163844
$ws.Range("S12").Value = 0.02012312207838108
$ws.Range("T12").Value = 0.02012312207838108

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 104.737245
$ws.Range("H13").Value = 314.211735
$ws.Range("I13").Value = 0.2898796499701289
$ws.Range("J13").Value = 0.2898796499701289
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 9.735567333333334
$ws.Range("N13").Value = 29.206702
$ws.Range("O13").Value = 0.2341158824922987
$ws.Range("P13").Value = 0.2341158824922987
$ws.Range("Q13").Value = 1019.67650100533
$ws.Range("R13").Value = 9177.08850904797
$ws.Range("S13").Value = 0.06786543006931539
$ws.Range("T13").Value = 0.06786543006931538

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 38.15794
$ws.Range("H14").Value = 114.47382
$ws.Range("I14").Value = 0.105609139239639
$ws.Range("J14").Value = 0.105609139239639
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 14.561928
$ws.Range("N14").Value = 43.685784
$ws.Range("O14").Value = 0.3501777048818433
$ws.Range("P14").Value = 0.3501777048818433
$ws.Range("Q14").Value = 555.65317490832
$ws.Range("R14").Value = 5000.87857417488
$ws.Range("S14").Value = 0.03698196599348381
$ws.Range("T14").Value = 0.03698196599348379

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 38.15794
$ws.Range("H15").Value = 114.47382
$ws.Range("I15").Value = 0.105609139239639
$ws.Range("J15").Value = 0.105609139239639
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 14.40015733333333
$ws.Range("N15").Value = 43.200472
$ws.Range("O15").Value = 0.3462875276490937
$ws.Range("P15").Value = 0.3462875276490937
$ws.Range("Q15").Value = 549.4803395158934
$ws.Range("R15").Value = 4945.323055643041
$ws.Range("S15").Value = 0.03657112772444348
$ws.Range("T15").Value = 0.03657112772444347

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 38.15794
$ws.Range("H16").Value = 114.47382
$ws.Range("I16").Value = 0.105609139239639
$ws.Range("J16").Value = 0.105609139239639
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 2.886742333333333
$ws.Range("N16").Value = 8.660226999999999
$ws.Range("O16").Value = 0.06941888497676431
$ws.Range("P16").Value = 0.06941888497676431
$ws.Range("Q16").Value = 110.1521407507933
$ws.Range("R16").Value = 991.3692667571399
$ws.Range("S16").Value = 0.007331268689371586
$ws.Range("T16").Value = 0.007331268689371584

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 38.15794
$ws.Range("H17").Value = 114.47382
$ws.Range("I17").Value = 0.105609139239639
$ws.Range("J17").Value = 0.105609139239639
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 9.735567333333334
$ws.Range("N17").Value = 29.206702
$ws.Range("O17").Value = 0.2341158824922987
$ws.Range("P17").Value = 0.2341158824922987
$ws.Range("Q17").Value = 371.4891941712934
$ws.Range("R17").Value = 371.4891941712934
$ws.Range("S17").Value = 0.02472477683234014
$ws.Range("T17").Value = 0.02472477683234013
